$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet / tab to reflect new "through" date
$ws.Name = "Through 2022-06-01"

# Update the "May" label (was "May (through 05-31)") now that May is complete
$ws.Range("A6").Value = "May"

# April (row 5) 2022 value correction: 116 -> 115
$ws.Range("I5").Value = 115

# Insert a new row for June just above the Total row, shifting Total down
$ws.Rows.Item(7).Insert()

# New June row data (through 2022-06-01); column B (2015) has no data
$ws.Range("A7").Value = "June (through 06-01)"
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = 6
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 10
$ws.Range("H7").Value = 2
$ws.Range("I7").Value = 4

# Match the bold/bordered/centered label style used by the other month cells in column A
$ws.Range("A7").Borders.LineStyle = 1        # xlContinuous
$ws.Range("A7").Font.Bold = $true
$ws.Range("A7").HorizontalAlignment = -4108  # xlCenter
$ws.Range("A7").VerticalAlignment = -4160    # xlTop

# Update the Total row (now row 8) with new totals including June data
$ws.Range("B8").Value = 108
$ws.Range("C8").Value = 210
$ws.Range("D8").Value = 318
$ws.Range("E8").Value = 301
$ws.Range("F8").Value = 205
$ws.Range("G8").Value = 368
$ws.Range("H8").Value = 633
$ws.Range("I8").Value = 668

# Widen column A slightly to fit the longer "June (through 06-01)" label
$ws.Columns.Item(1).ColumnWidth = $ws.Columns.Item(1).ColumnWidth + 1
